$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the SUPRAX 100MG/5 ML SUSP. 60ML row (row 99) entirely; this
#    shifts every row below it up by one (data rows, totals row, footer row).
$ws.Rows("99").Delete()

# 2. Column A holds a simple running serial number (1..N) for the data
#    rows and is regenerated after the deletion, so re-sequence it rather
#    than letting it inherit the row below's number.
for ($i = 0; $i -lt 122; $i++) {
    $ws.Cells.Item(7 + $i, 1).Value = $i + 1
}

# 3. Correct the "WATER FOR INJECTION AMP. 5 ML" row (now row 107): the
#    stock count, the sold amount and the transaction count were updated.
#    These columns store text (e.g. "14.0000"), not numbers, so force a
#    text number-format while assigning and restore it afterwards to avoid
#    Excel auto-converting the string into a real number.
$fmtH = $ws.Range("H107").NumberFormat
$ws.Range("H107").NumberFormat = "@"
$ws.Range("H107").Value = "8492:0"
$ws.Range("H107").NumberFormat = $fmtH

$fmtP = $ws.Range("P107").NumberFormat
$ws.Range("P107").NumberFormat = "@"
$ws.Range("P107").Value = "14.0000"
$ws.Range("P107").NumberFormat = $fmtP

$fmtQ = $ws.Range("Q107").NumberFormat
$ws.Range("Q107").NumberFormat = "@"
$ws.Range("Q107").Value = "7:0"
$ws.Range("Q107").NumberFormat = $fmtQ

# 4. Update the grand-total cell (now row 129).
$ws.Range("P129").Value = 7959.2299999999996

# 5. Update the generated-on timestamp in the footer (now row 130).
$ws.Range("A130").Value = "Monday, 28 July, 2025 10:48 PM"

# 6. Row heights were auto-fit per row content and, once the row shift
#    happens, Excel recomputes them based on the new text in each row;
#    restore every data-row height (99-129) to what it was at that same
#    row position before the edit, matching the regenerated report.
$rowHeights = @{
    99 = 25.5; 100 = 24.75; 101 = 25.5; 102 = 25.5; 103 = 24.75; 104 = 25.5;
    105 = 24.75; 106 = 25.5; 107 = 25.5; 108 = 24.75; 109 = 25.5; 110 = 24.75;
    111 = 25.5; 112 = 25.5; 113 = 24.75; 114 = 25.5; 115 = 24.75; 116 = 25.5;
    117 = 25.5; 118 = 24.75; 119 = 25.5; 120 = 24.75; 121 = 25.5; 122 = 25.5;
    123 = 24.75; 124 = 25.5; 125 = 24.75; 126 = 25.5; 127 = 25.5; 128 = 24.75;
    129 = 25.5; 130 = 16.5
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows($r).RowHeight = $rowHeights[$r]
}
